$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column stays text (matches original "inlineStr" formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.859.10"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.741.26"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "238.19"
$ws.Range("E5").Value = "  +3.67%  "
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.5150"
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").Value = "0.2740"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "40.00"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").Value = "0.06130"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "1.740.65"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "0.07175"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "15.03"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "0.6429"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "4.596"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "77.36"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").Value = "0.9978"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "0.9986"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "25.898.01"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "11.76"
$ws.Range("E20").Value = "  +2.92%  "
$ws.Range("D21").Value = "0.000006774"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").Value = "1.960.89"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "4.275"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").Value = "8.669"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").Value = "5.244"
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").Value = "138.74"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "1.528"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").Value = "15.23"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").Value = "1.769"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "106.36"
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("E31").Value = "  +8.91%  "
$ws.Range("D32").Value = "0.08326"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "3.652"
$ws.Range("E33").Value = "  +4.37%  "
$ws.Range("D34").Value = "0.04592"
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("D35").Value = "2.658"
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("D36").Value = "0.9907"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("D37").Value = "0.6206"
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("D38").Value = "2.693"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").Value = "0.01618"
$ws.Range("E39").Value = "  +3.36%  "
$ws.Range("D40").Value = "1.937"
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("D41").Value = "0.9974"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "97.96"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").Value = "0.3854"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "0.7385"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").Value = "4.948"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").Value = "0.1126"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "0.05262"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "6.191"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "55.00"
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("D50").Value = "30.54"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("D51").Value = "7.590"
$ws.Range("E51").Value = "  -0.07%  "
